# Update crypto price/volume data per the Aug 8 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = <new price text>; E = <new volume text> }  (omitted keys are left untouched)
$updates = @{
    2 = @{ D = "29.760.87"; E = "  +2.57%  " }
    3 = @{ D = "1.853.37"; E = "  +1.81%  " }
    4 = @{ D = "0.9988"; E = "  -0.06%  " }
    5 = @{ D = "244.79"; E = "  +1.34%  " }
    6 = @{ D = "0.6386" }
    7 = @{ D = "0.9992"; E = "  -0.36%  " }
    8 = @{ D = "0.07515"; E = "  +2.78%  " }
    9 = @{ D = "0.2941"; E = "  +2.29%  " }
    10 = @{ D = "24.05"; E = "  +5.05%  " }
    11 = @{ D = "0.07677"; E = "  +0.35%  " }
    12 = @{ D = "1.839.36" }
    13 = @{ D = "5.059"; E = "  +2.64%  " }
    14 = @{ D = "0.6888"; E = "  +4.64%  " }
    15 = @{ E = "  +4.51%  " }
    16 = @{ D = "0.000009622"; E = "  +8.08%  " }
    17 = @{ D = "6.054"; E = "  +3.56%  " }
    18 = @{ D = "29.714.30"; E = "  +2.52%  " }
    19 = @{ D = "2.088.24"; E = "  +1.18%  " }
    20 = @{ D = "239.54"; E = "  +1.36%  " }
    21 = @{ D = "12.66"; E = "  +1.89%  " }
    22 = @{ D = "0.9994"; E = "  -0.31%  " }
    23 = @{ D = "7.378"; E = "  +3.81%  " }
    24 = @{ D = "0.9993"; E = "  -0.25%  " }
    25 = @{ D = "159.53"; E = "  +0.52%  " }
    26 = @{ D = "0.1430"; E = "  +2.25%  " }
    27 = @{ D = "8.558"; E = "  +1.83%  " }
    28 = @{ D = "17.96"; E = "  +2.13%  " }
    29 = @{ D = "1.504"; E = "  +1.75%  " }
    30 = @{ D = "0.06010"; E = "  +7.56%  " }
    31 = @{ D = "1.257"; E = "  +3.81%  " }
    32 = @{ D = "4.160"; E = "  +2.40%  " }
    33 = @{ D = "4.146"; E = "  +1.81%  " }
    34 = @{ D = "1.875"; E = "  +3.08%  " }
    35 = @{ D = "1.152"; E = "  +1.98%  " }
    36 = @{ D = "0.7340"; E = "  +0.54%  " }
    37 = @{ D = "2.606"; E = "  -0.92%  " }
    38 = @{ D = "2.870"; E = "  +2.19%  " }
    39 = @{ D = "1.228.90"; E = "  +2.89%  " }
    40 = @{ E = "  +1.60%  " }
    41 = @{ D = "6.387"; E = "  +0.72%  " }
    42 = @{ D = "0.9158"; E = "  +3.67%  " }
    43 = @{ D = "1.000"; E = "  -0.28%  " }
    44 = @{ D = "2.001.58"; E = "  +1.87%  " }
    45 = @{ D = "102.02"; E = "  +1.42%  " }
    46 = @{ E = "  +3.36%  " }
    47 = @{ E = "  +2.01%  " }
    48 = @{ D = "0.5078"; E = "  -0.15%  " }
    49 = @{ D = "9.341" }
    50 = @{ D = "0.4091"; E = "  +2.82%  " }
    51 = @{ D = "0.1140"; E = "  +4.09%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        # Force text so Excel keeps values like "0.9988" / "1.000" / "0.06010" literally
        # instead of re-parsing them as numbers and dropping trailing/format digits.
        $ws.Range("D$row").NumberFormat = "@"
        $ws.Range("D$row").Value = $vals["D"]
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}

# Drop the explicit text-number-format style again so cells keep the workbook default
# (no "s" attribute), matching how the rest of the data column is styled.
$ws.Range("D2:D51").Style = "Normal"
